# Generate Report for Handback
# Removes the handback row for "e0390748-aa6c-4338-96c4-c0d594e4314f" (row 3)
# from all three sheets (Overview, zh-cn, de-de) and refreshes the
# "Correspond Handoff/Handback Datetime" timestamps on the remaining
# (8e2b4336...) row for the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$hlsOverview = @($wsOverview.Hyperlinks)
for ($i = $hlsOverview.Count - 1; $i -ge 0; $i--) {
    $hl = $hlsOverview[$i]
    if ($hl.Range.Row -eq 3) {
        $hl.Delete()
    }
}
$wsOverview.Rows(3).Delete()

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Refresh the handoff/handback timestamps for the remaining row (row 2).
$wsZhCn.Range("E2").Value = "2016-03-14 04:31:40"
$wsZhCn.Range("H2").Value = "2016-03-14 04:31:56"

# Drop the row for the removed file (row 3) and its hyperlinks.
# (Iterate in reverse so deleting a hyperlink doesn't shift the indices
# of the not-yet-processed items in the snapshot.)
$hlsZhCn = @($wsZhCn.Hyperlinks)
for ($i = $hlsZhCn.Count - 1; $i -ge 0; $i--) {
    $hl = $hlsZhCn[$i]
    if ($hl.Range.Row -eq 3) {
        $hl.Delete()
    }
}
$wsZhCn.Rows(3).Delete()

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Refresh the handoff/handback timestamps for the remaining row (row 2).
$wsDeDe.Range("E2").Value = "2016-03-14 04:31:43"
$wsDeDe.Range("H2").Value = "2016-03-14 04:32:02"

# Drop the row for the removed file (row 3) and its hyperlinks.
# (Iterate in reverse so deleting a hyperlink doesn't shift the indices
# of the not-yet-processed items in the snapshot.)
$hlsDeDe = @($wsDeDe.Hyperlinks)
for ($i = $hlsDeDe.Count - 1; $i -ge 0; $i--) {
    $hl = $hlsDeDe[$i]
    if ($hl.Range.Row -eq 3) {
        $hl.Delete()
    }
}
$wsDeDe.Rows(3).Delete()
